$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new values look numeric as Text, so the literal
# string (e.g. "553.03") is preserved exactly instead of being parsed into
# a floating point number by Excel.
$textCells = @("D4", "D5", "D6", "D7", "D10", "D11", "D12", "D13", "D14", "D17", "D19", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D37", "D38", "D39", "D40", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.381.11"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "3.335.67"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "553.03"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").Value = "173.11"
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  +1.43%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "3.325.92"
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +6.62%  "
$ws.Range("D11").Value = "0.637"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").Value = "53.37"
$ws.Range("E12").Value = "  -3.10%  "
$ws.Range("D13").Value = "0.0000277"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").Value = "9.08"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").Value = "3.851.56"
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "18.12"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "3.318.10"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "11.73"
$ws.Range("E19").Value = "  -1.16%  "
$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").Value = "64.147.35"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").Value = "0.986"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").Value = "448.38"
$ws.Range("E22").Value = "  +3.60%  "
$ws.Range("D23").Value = "4.99"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("D24").Value = "4.04"
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").Value = "86.68"
$ws.Range("E25").Value = "  +2.74%  "
$ws.Range("E26").Value = "  +3.44%  "
$ws.Range("D27").Value = "2.87"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("D28").Value = "10.63"
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("D29").Value = "8.59"
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("D30").Value = "30.85"
$ws.Range("E30").Value = "  +3.63%  "
$ws.Range("D31").Value = "6.52"
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "11.38"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").Value = "62.35"
$ws.Range("E33").Value = "  +6.50%  "
$ws.Range("D34").Value = "568.78"
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.141"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "3.53"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").Value = "35.26"
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("D40").Value = "0.367"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("D41").Value = "0.0₃0730"
$ws.Range("E41").Value = "  -3.98%  "
$ws.Range("D42").Value = "3.059.42"
$ws.Range("E42").Value = "  -1.92%  "
$ws.Range("D43").Value = "0.0414"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").Value = "2.73"
$ws.Range("E44").Value = "  -3.86%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "0.133"
$ws.Range("E45").Value = "  +2.84%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "2.44"
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "3.16"
$ws.Range("E47").Value = "  -3.66%  "
$ws.Range("D48").Value = "0.996"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").Value = "140.35"
$ws.Range("E49").Value = "  +3.87%  "
$ws.Range("E50").Value = "  -3.34%  "
$ws.Range("D51").Value = "8.18"
$ws.Range("E51").Value = "  -1.40%  "
